# The commit swaps the two theme parts in the package: the slide-master's
# theme (ppt/theme/theme1.xml, originally the "Integral" / "Red Violet"
# palette) is replaced by the content that used to live in the notes
# master's theme (ppt/theme/theme2.xml, the stock "Office Theme" / "Office"
# palette) - and vice versa.
#
# The only part of that swap reachable through the PowerPoint object model
# is the document theme's (SlideMaster's) ThemeColorScheme - font scheme and
# format scheme are already identical between the two theme parts, so the
# only real difference to reproduce is the 12-slot colour scheme (and, were
# it possible, the cosmetic theme/colour-scheme display names, which
# PowerPoint does not expose as a settable property).
#
# Theme color slot order (ThemeColorScheme.Item index -> OOXML tag):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
#
# RGB() COM values are stored as 0x00BBGGRR, so each hex colour below is
# converted to B<<16 | G<<8 | R before being assigned.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
